# Update event statuses:
# - Arm Wrestling (row 4): move from UPCOMING to ONGOING (F4: 1 -> 0, G4: 0 -> 1)
# - Chess (row 8): TAG changes from "x" to "i" (indoor)
# - Update the active selection to I9 (cosmetic, matches last-saved cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Arm Wrestling: no longer upcoming, now ongoing
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1

# Row 8 - Chess: fix tag from "x" to "i" (indoor)
$ws.Range("I8").Value = "i"

# Update the selected cell to match the saved cursor position
$ws.Range("I9").Select()
